$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("5 Era Vandervoort")
$ws.Activate()

# Update the "Expertise" (column G) ratings for several projects.
# Cells that gain a new value (previously blank):
$ws.Range("G2").Value = "L"
$ws.Range("G5").Value = "L"
$ws.Range("G12").Value = "H"
$ws.Range("G16").Value = "L"
$ws.Range("G23").Value = "M"
$ws.Range("G24").Value = "L"
$ws.Range("G25").Value = "L"

# Cells that change value:
$ws.Range("G20").Value = "L"

# Cells that are cleared (previously had a value):
$ws.Range("G9").Value = ""
$ws.Range("G11").Value = ""
$ws.Range("G14").Value = ""
$ws.Range("G19").Value = ""
$ws.Range("G22").Value = ""

# Update the selection shown on the sheet.
$ws.Range("A10:K10").Select()
